# Update the battery cost and notebook
#
# The "default_parameters_values" sheet lists battery-model parameters
# (columns: name | amount | formula). This change adds five new separator-
# thickness parameters (11/13/15/17/19 um), each defaulting to an amount of
# 0, and folds them into the "separator" row's aggregation formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new parameters belong right before "separator_waste" (row 124), which
# keeps the separator_* family of rows grouped together. Insert 5 blank rows
# there; everything from the old row 124 down (including the "separator"
# row itself) shifts down by 5.
$ws.Rows("124:128").Insert()

$ws.Range("A124").Value = "separator_11um"
$ws.Range("B124").Value = 0
$ws.Range("A125").Value = "separator_13um"
$ws.Range("B125").Value = 0
$ws.Range("A126").Value = "separator_15um"
$ws.Range("B126").Value = 0
$ws.Range("A127").Value = "separator_17um"
$ws.Range("B127").Value = 0
$ws.Range("A128").Value = "separator_19um"
$ws.Range("B128").Value = 0

# The "separator" row (now row 134) aggregates the separator thickness
# variants plus the coated-separator variants — extend its formula text to
# include the five newly added terms.
$ws.Range("C134").Value = "(separator_5um+separator_7um+separator_9um+separator_11um+separator_13um+separator_15um+separator_17um+separator_19um+coated_separator_5um_2um+coated_separator_7um_2um+coated_separator_9um_3um)"

# Match the author's on-save cursor position.
$ws.Range("C120").Select()
